# Update gh-pages to output generated at 456a3b4
# Updates the "想去人数" (want-to-go count) values in column F on the
# "展览" sheet (rows 4,5,6,10,11) and the "全部类型" sheet (rows 4,5,6,10,14).

$wb = $excel.ActiveWorkbook

# --- Sheet "展览" ---
$wsExhibit = $wb.Worksheets.Item("展览")
$wsExhibit.Range("F4").Value = 8000
$wsExhibit.Range("F5").Value = 5838
$wsExhibit.Range("F6").Value = 494
$wsExhibit.Range("F10").Value = 283
$wsExhibit.Range("F11").Value = 364

# --- Sheet "全部类型" ---
$wsAll = $wb.Worksheets.Item("全部类型")
$wsAll.Range("F4").Value = 8000
$wsAll.Range("F5").Value = 5838
$wsAll.Range("F6").Value = 494
$wsAll.Range("F10").Value = 283
$wsAll.Range("F14").Value = 364
